$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.703.21"
$ws.Range("D2").Style = $dStyle
$ws.Range("E2").Value = "  +0.56%  "

$dStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.742.88"
$ws.Range("D3").Style = $dStyle
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  -0.11%  "

$dStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.45"
$ws.Range("D5").Style = $dStyle
$ws.Range("E5").Value = "  +0.80%  "

$dStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.98"
$ws.Range("D6").Style = $dStyle
$ws.Range("E6").Value = "  +2.38%  "

$dStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.741.51"
$ws.Range("D7").Style = $dStyle
$ws.Range("E7").Value = "  +1.15%  "

$ws.Range("E8").Value = "  -0.08%  "

$dStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = $dStyle
$ws.Range("E9").Value = "  -1.12%  "

$dStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.166"
$ws.Range("D10").Style = $dStyle
$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("E11").Value = "  +3.61%  "

$dStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.482"
$ws.Range("D12").Style = $dStyle
$ws.Range("E12").Value = "  -2.59%  "

$dStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.95"
$ws.Range("D13").Style = $dStyle
$ws.Range("E13").Value = "  -0.98%  "

$dStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000254"
$ws.Range("D14").Style = $dStyle
$ws.Range("E14").Value = "  +1.07%  "

$dStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.358.83"
$ws.Range("D15").Style = $dStyle
$ws.Range("E15").Value = "  +1.12%  "

$dStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.741.81"
$ws.Range("D16").Style = $dStyle
$ws.Range("E16").Value = "  +1.11%  "

$dStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.742.12"
$ws.Range("D17").Style = $dStyle
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("E18").Value = "  -2.30%  "

$dStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.47"
$ws.Range("D19").Style = $dStyle
$ws.Range("E19").Value = "  -0.67%  "

$dStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "502.40"
$ws.Range("D20").Style = $dStyle
$ws.Range("E20").Value = "  -1.65%  "

$dStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.35"
$ws.Range("D21").Style = $dStyle
$ws.Range("E21").Value = "  -1.33%  "

$dStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.15"
$ws.Range("D22").Style = $dStyle
$ws.Range("E22").Value = "  -1.84%  "

$dStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.721"
$ws.Range("D23").Style = $dStyle
$ws.Range("E23").Value = "  +0.02%  "

$dStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.66"
$ws.Range("D24").Style = $dStyle
$ws.Range("E24").Value = "  +8.56%  "

$dStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.07"
$ws.Range("D25").Style = $dStyle
$ws.Range("E25").Value = "  -1.49%  "

$dStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.66"
$ws.Range("D26").Style = $dStyle
$ws.Range("E26").Value = "  +7.29%  "

$dStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.91"
$ws.Range("D27").Style = $dStyle
$ws.Range("E27").Value = "  -2.91%  "

$dStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000137"
$ws.Range("D28").Style = $dStyle
$ws.Range("E28").Value = "  +10.38%  "

$ws.Range("E29").Value = "  +0.59%  "

$dStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.48"
$ws.Range("D30").Style = $dStyle
$ws.Range("E30").Value = "  -0.33%  "

$ws.Range("E31").Value = "  +3.69%  "

$dStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.05"
$ws.Range("D32").Style = $dStyle
$ws.Range("E32").Value = "  +3.35%  "

$dStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.40"
$ws.Range("D33").Style = $dStyle
$ws.Range("E33").Value = "  -1.84%  "

$dStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.114"
$ws.Range("D34").Style = $dStyle
$ws.Range("E34").Value = "  -0.55%  "

$dStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("D35").Style = $dStyle
$ws.Range("E35").Value = "  -0.33%  "

$ws.Range("E36").Value = "  +1.57%  "

$dStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.12"
$ws.Range("D37").Style = $dStyle
$ws.Range("E37").Value = "  -0.20%  "

$dStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.355"
$ws.Range("D38").Style = $dStyle
$ws.Range("E38").Value = "  +5.61%  "

$ws.Range("E39").Value = "  +4.14%  "

$dStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("D40").Style = $dStyle
$ws.Range("E40").Value = "  +14.45%  "

$dStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.07"
$ws.Range("D41").Style = $dStyle
$ws.Range("E41").Value = "  -3.95%  "

$dStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "46.26"
$ws.Range("D42").Style = $dStyle
$ws.Range("E42").Value = "  +5.60%  "

$dStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "443.70"
$ws.Range("D43").Style = $dStyle
$ws.Range("E43").Value = "  +7.37%  "

$dStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "49.73"
$ws.Range("D44").Style = $dStyle
$ws.Range("E44").Value = "  -2.87%  "

$dStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.56"
$ws.Range("D45").Style = $dStyle
$ws.Range("E45").Value = "  -2.05%  "

$dStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.955.14"
$ws.Range("D46").Style = $dStyle
$ws.Range("E46").Value = "  -3.58%  "

$dStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0360"
$ws.Range("D47").Style = $dStyle
$ws.Range("E47").Value = "  -0.17%  "

$dStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.65"
$ws.Range("D48").Style = $dStyle
$ws.Range("E48").Value = "  +2.77%  "

$dStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.11"
$ws.Range("D50").Style = $dStyle
$ws.Range("E50").Value = "  -1.77%  "

$dStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.50"
$ws.Range("D51").Style = $dStyle
$ws.Range("E51").Value = "  +0.34%  "
